# Generate Report for Archive
#
# The archive-report generator re-ran and this time processed
# "ea4cc21f-9abe-4720-8507-1b331f119a75.md" before
# "4b647b34-46ab-454f-8905-9f77375c347d.md", so the two data rows (6 and 7)
# describing those files swap places on every sheet (Overview, zh-cn,
# de-de). Only the cells that actually change value are touched; empty
# cells, headers, table ranges and hyperlink targets (rIds) are left as-is.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview": File Name / Path And Name / zh-cn / de-de / Latest HO Xliff Generate Date ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A6").Value = "ea4cc21f-9abe-4720-8507-1b331f119a75.md"
$wsOverview.Range("B6").Value = "e2e\ea4cc21f-9abe-4720-8507-1b331f119a75.md"
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"
$wsOverview.Range("G6").Value = "2016-10-19 23:28:26"

$wsOverview.Range("A7").Value = "4b647b34-46ab-454f-8905-9f77375c347d.md"
$wsOverview.Range("B7").Value = "e2e\4b647b34-46ab-454f-8905-9f77375c347d.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-10-19 23:23:54"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$6') { $hl.TextToDisplay = "e2e\ea4cc21f-9abe-4720-8507-1b331f119a75.md" }
    if ($hl.Range.Address() -eq '$B$7') { $hl.TextToDisplay = "e2e\4b647b34-46ab-454f-8905-9f77375c347d.md" }
}

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A6").Value = "ea4cc21f-9abe-4720-8507-1b331f119a75.md"
$wsZh.Range("C6").Value = "In Translation"
$wsZh.Range("G6").Value = "ea4cc21f-9abe-4720-8507-1b331f119a75.937dcae6e5f8cde315e31f0300e387fe35b13d7e.zh-cn.xlf"
$wsZh.Range("H6").Value = "2016-10-19 23:28:14"

$wsZh.Range("A7").Value = "4b647b34-46ab-454f-8905-9f77375c347d.md"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("G7").Value = "4b647b34-46ab-454f-8905-9f77375c347d.8faef7f17390d39a282eec5c85ef893ba9b23988.zh-cn.xlf"
$wsZh.Range("H7").Value = "2016-10-19 23:23:44"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$6') { $hl.TextToDisplay = "ea4cc21f-9abe-4720-8507-1b331f119a75.md" }
    if ($hl.Range.Address() -eq '$A$7') { $hl.TextToDisplay = "4b647b34-46ab-454f-8905-9f77375c347d.md" }
}

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A6").Value = "ea4cc21f-9abe-4720-8507-1b331f119a75.md"
$wsDe.Range("C6").Value = "In Translation"
$wsDe.Range("G6").Value = "ea4cc21f-9abe-4720-8507-1b331f119a75.937dcae6e5f8cde315e31f0300e387fe35b13d7e.de-de.xlf"
$wsDe.Range("H6").Value = "2016-10-19 23:28:26"

$wsDe.Range("A7").Value = "4b647b34-46ab-454f-8905-9f77375c347d.md"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("G7").Value = "4b647b34-46ab-454f-8905-9f77375c347d.8faef7f17390d39a282eec5c85ef893ba9b23988.de-de.xlf"
$wsDe.Range("H7").Value = "2016-10-19 23:23:54"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$6') { $hl.TextToDisplay = "ea4cc21f-9abe-4720-8507-1b331f119a75.md" }
    if ($hl.Range.Address() -eq '$A$7') { $hl.TextToDisplay = "4b647b34-46ab-454f-8905-9f77375c347d.md" }
}
